$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.275026679039001
$ws.Range("B1").Value = 2.349532604217529
$ws.Range("D1").Value = 1.384984850883484
$ws.Range("E1").Value = 0.8522005677223206
